$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = 1987..2024

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cell.Value = Get-Date -Year $years[$i] -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
}

$wb.Save()
